# Fix removing redundant entries:
# When a buy/sell day pair appears with duplicate adjacent entries, keep the
# correct single record. This rewrites the profit table (rows 2-21) with the
# corrected buy day / sell day / profit values, re-sorted by profit descending.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "friday"
$ws.Cells.Item(2, 2).Value = "thursday"
$ws.Cells.Item(2, 3).Value = 188.1070723521612
$ws.Cells.Item(3, 1).Value = "friday"
$ws.Cells.Item(3, 2).Value = "wednesday"
$ws.Cells.Item(3, 3).Value = 156.4264313298883
$ws.Cells.Item(4, 1).Value = "saturday"
$ws.Cells.Item(4, 2).Value = "thursday"
$ws.Cells.Item(4, 3).Value = 129.0821741183443
$ws.Cells.Item(5, 1).Value = "thursday"
$ws.Cells.Item(5, 2).Value = "wednesday"
$ws.Cells.Item(5, 3).Value = 120.6760754743007
$ws.Cells.Item(6, 1).Value = "friday"
$ws.Cells.Item(6, 2).Value = "tuesday"
$ws.Cells.Item(6, 3).Value = 101.7144789607758
$ws.Cells.Item(7, 1).Value = "saturday"
$ws.Cells.Item(7, 2).Value = "wednesday"
$ws.Cells.Item(7, 3).Value = 94.61486169571577
$ws.Cells.Item(8, 1).Value = "saturday"
$ws.Cells.Item(8, 2).Value = "friday"
$ws.Cells.Item(8, 3).Value = 92.23404854491551
$ws.Cells.Item(9, 1).Value = "wednesday"
$ws.Cells.Item(9, 2).Value = "tuesday"
$ws.Cells.Item(9, 3).Value = 88.49543971791351
$ws.Cells.Item(10, 1).Value = "thursday"
$ws.Cells.Item(10, 2).Value = "tuesday"
$ws.Cells.Item(10, 3).Value = 69.90923716236374
$ws.Cells.Item(11, 1).Value = "saturday"
$ws.Cells.Item(11, 2).Value = "tuesday"
$ws.Cells.Item(11, 3).Value = 51.71913448807164
$ws.Cells.Item(12, 1).Value = "tuesday"
$ws.Cells.Item(12, 2).Value = "thursday"
$ws.Cells.Item(12, 3).Value = 44.52487332378689
$ws.Cells.Item(13, 1).Value = "friday"
$ws.Cells.Item(13, 2).Value = "saturday"
$ws.Cells.Item(13, 3).Value = 29.43699693634931
$ws.Cells.Item(14, 1).Value = "tuesday"
$ws.Cells.Item(14, 2).Value = "wednesday"
$ws.Cells.Item(14, 3).Value = 20.45246752243328
$ws.Cells.Item(15, 1).Value = "tuesday"
$ws.Cells.Item(15, 2).Value = "saturday"
$ws.Cells.Item(15, 3).Value = 13.16876223331853
$ws.Cells.Item(16, 1).Value = "wednesday"
$ws.Cells.Item(16, 2).Value = "thursday"
$ws.Cells.Item(16, 3).Value = 9.328489758797152
$ws.Cells.Item(17, 1).Value = "tuesday"
$ws.Cells.Item(17, 2).Value = "friday"
$ws.Cells.Item(17, 3).Value = -3.223615634411131
$ws.Cells.Item(18, 1).Value = "wednesday"
$ws.Cells.Item(18, 2).Value = "saturday"
$ws.Cells.Item(18, 3).Value = -7.786797665645799
$ws.Cells.Item(19, 1).Value = "thursday"
$ws.Cells.Item(19, 2).Value = "saturday"
$ws.Cells.Item(19, 3).Value = -16.24622036142086
$ws.Cells.Item(20, 1).Value = "wednesday"
$ws.Cells.Item(20, 2).Value = "friday"
$ws.Cells.Item(20, 3).Value = -29.4483929612289
$ws.Cells.Item(21, 1).Value = "thursday"
$ws.Cells.Item(21, 2).Value = "friday"
$ws.Cells.Item(21, 3).Value = -36.99842750196777
